# The workbook's data rows (2-18, row 16 untouched by the underlying edit)
# got their "record" columns (Fecha, Variedad, Calidad, Volumen, Precio
# minimo/maximo/promedio, Unidad de comercializacion, Precio $/Kg, Kg o
# Unidades) rotated among rows - i.e. each destination row's record came
# from a different source row (row 16 keeps its own data). Columns A, B,
# C, E, F, G, O, R are identical across all rows so they are left alone.
#
# destination row -> source row the "record" came from
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2  = 17
    3  = 2
    4  = 8
    5  = 9
    6  = 10
    7  = 3
    8  = 4
    9  = 6
    10 = 14
    11 = 13
    12 = 11
    13 = 5
    14 = 18
    15 = 12
    16 = 16
    17 = 15
    18 = 7
}

# Columns (by index) that move together as one "record" per the diff.
# D=4 H=8 I=9 J=10 K=11 L=12 M=13 N=14 P=16 Q=17
$cols = @(4, 8, 9, 10, 11, 12, 13, 14, 16, 17)

# Snapshot the "before" record for every source row first, so overwriting
# a row that is itself used as a source for another row is safe.
$snapshots = @{}
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    if (-not $snapshots.ContainsKey($srcRow)) {
        $rec = @{}
        foreach ($c in $cols) {
            $rec[$c] = $ws.Cells.Item($srcRow, $c).Value2
        }
        $snapshots[$srcRow] = $rec
    }
}

# Now write each destination row's record from its captured snapshot.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $rec = $snapshots[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $rec[$c]
    }
}
